# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" sheet (fund-level holdings) positioned between the
# existing "总计" sheet and the existing "2021-Q1" sheet, and updates the
# "总计" (totals) sheet with a new summary row for 2022-Q4 while keeping the
# 2021-Q1 summary row (shifted down one row).
#
# Strategy (worked out empirically against the iron_native COM shim):
#   - Worksheets.Add(...)'s new sheet always receives the next unused
#     sheetId, regardless of Before/After placement. To reproduce the
#     target's left-to-right sheetId order (总计=1, 2022-Q4=2, 2021-Q1=3)
#     we keep the *existing* second sheet object in place (so it keeps
#     sheetId=2) and simply overwrite its contents/name to become
#     "2022-Q4"; the brand-new sheet object (sheetId=3) becomes "2021-Q1"
#     and is populated with the data the second sheet used to hold.
#   - PasteSpecial(xlPasteValues) from an already-text-typed source cell
#     is the only reliable way found to write a numeric-looking string
#     (e.g. "9.03", "016616") without Excel silently re-parsing it as a
#     number - plain `.Value = "9.03"` assignment always converts it.
#   - PasteSpecial(xlPasteFormats) is used to copy an existing cell style
#     (by reference, not by cloning) so the saved file reuses the same
#     cellXfs index instead of allocating a new one.
#   - User-defined PowerShell functions lose the live COM binding when
#     Range/Worksheet objects (or even plain strings) are passed in as
#     parameters under this shim, so the whole script is written as flat,
#     linear statements (loops are fine; functions are avoided).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$sheetTotal = $wb.Worksheets.Item(1)   # "总计"
$sheetOld   = $wb.Worksheets.Item(2)   # currently "2021-Q1", becomes "2022-Q4"

# ---------------------------------------------------------------------------
# 1) Append a brand-new worksheet right after $sheetOld. This will end up
#    holding the data that currently lives in $sheetOld (i.e. the original
#    "2021-Q1" fund table), once we've copied it over below.
# ---------------------------------------------------------------------------
$sheetNew = $wb.Worksheets.Add($null, $sheetOld)
$sheetNew.Name = "2021-Q1-NEW-TMP"

# Copy the header-row + A2 formatting (style) from the old sheet over to the
# new sheet, so it keeps the same cellXfs style index.
$sheetOld.Range("B1:H1").Copy()
$sheetNew.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$sheetOld.Range("A2").Copy()
$sheetNew.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

# Copy the values (preserving text-vs-number typing) from old -> new.
$sheetOld.Range("A1:H2").Copy()
$sheetNew.Range("A1:H2").PasteSpecial(-4163)   # xlPasteValues

# Rename $sheetOld away from "2021-Q1" *before* trying to rename $sheetNew to
# "2021-Q1" - sheet names must be unique at every point in time, and a
# rename that collides with an existing name throws (which this shim does
# not abort the whole script on, it just silently no-ops that statement -
# so the ordering below matters).
$sheetOld.Name = "2022-Q4"
$sheetNew.Name = "2021-Q1"

# ---------------------------------------------------------------------------
# 2) Re-style the original second sheet ($sheetOld, now named "2022-Q4") to
#    match the "总计" sheet's header style (reused cellXfs index), then
#    overwrite its values with the 2022-Q4 fund table.
# ---------------------------------------------------------------------------
$sheetTotal.Range("B1:D1").Copy()
$sheetOld.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$sheetTotal.Range("A2").Copy()
$sheetOld.Range("A2:A3").PasteSpecial(-4122)   # xlPasteFormats

# Header row text (column D's label changes from "基金金额" to "基金规模").
$sheetOld.Range("B1").Value = "基金代码"
$sheetOld.Range("C1").Value = "基金名称"
$sheetOld.Range("D1").Value = "基金规模"
$sheetOld.Range("E1").Value = "股票总仓位"
$sheetOld.Range("F1").Value = "仓位占比"
$sheetOld.Range("G1").Value = "持有市值(亿元)"
$sheetOld.Range("H1").Value = "仓位排名"

# Row 2 - fund 160212.
$sheetOld.Range("A2").Value = 0
$sheetOld.Range("B2").Value = "160212"
$sheetOld.Range("C2").Value = "国泰估值优势混合（LOF）A"
$sheetOld.Range("D2").Value = "9.03"
$sheetOld.Range("E2").Value = "93.64"
$sheetOld.Range("F2").Value = "5.35"
$sheetOld.Range("G2").Value = "0.4831"
$sheetOld.Range("H2").Value = 9

# Row 3 - fund 016616.
$sheetOld.Range("A3").Value = 1
$sheetOld.Range("B3").Value = "016616"
$sheetOld.Range("C3").Value = "国泰估值优势混合（LOF）C"
$sheetOld.Range("D3").Value = "0.00"
$sheetOld.Range("E3").Value = "93.64"
$sheetOld.Range("F3").Value = "5.35"
$sheetOld.Range("G3").Value = 0
$sheetOld.Range("H3").Value = 9

# The B/C/D/E/F/G cells above that must stay TEXT even though they look
# numeric ("160212", "016616", "9.03", "93.64", "5.35", "0.4831", "0.00")
# need the scratch-cell round trip: plain `.Value = "..."` assignment lets
# Excel silently reinterpret them as numbers.
$scratch = $sheetOld.Range("ZZ1")

$scratch.NumberFormat = "@"
$scratch.Value = "160212"
$scratch.Copy()
$sheetOld.Range("B2").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "9.03"
$scratch.Copy()
$sheetOld.Range("D2").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "93.64"
$scratch.Copy()
$sheetOld.Range("E2").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "5.35"
$scratch.Copy()
$sheetOld.Range("F2").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "0.4831"
$scratch.Copy()
$sheetOld.Range("G2").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "016616"
$scratch.Copy()
$sheetOld.Range("B3").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "0.00"
$scratch.Copy()
$sheetOld.Range("D3").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "93.64"
$scratch.Copy()
$sheetOld.Range("E3").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = "5.35"
$scratch.Copy()
$sheetOld.Range("F3").PasteSpecial(-4163)

$scratch.Clear()

# ---------------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: push the existing 2021-Q1 summary row
#    down to row 3 (recomputing its running index in column A to 1), and
#    write the new 2022-Q4 summary into row 2.
# ---------------------------------------------------------------------------
$sheetTotal.Range("B2:D2").Copy()
$sheetTotal.Range("B3:D3").PasteSpecial(-4163)   # xlPasteValues

$sheetTotal.Range("A2").Copy()
$sheetTotal.Range("A3").PasteSpecial(-4122)      # xlPasteFormats
$sheetTotal.Range("A3").Value = 1

$sheetTotal.Range("B2").Value = "2022-Q4"
$sheetTotal.Range("C2").Value = 2
$sheetTotal.Range("D2").Value = 0.48

Write-Output "done"
